$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "TT" (sequence number) column entirely; everything to the
# right shifts one column to the left (data, styles, shared strings, and
# column width definitions all move with it).
$ws.Columns("A").Delete()

# The newly freed-up columns AF:AJ (previously AG:AK) get explicit custom
# widths.
$ws.Columns("AF").ColumnWidth = 21.5703125
$ws.Columns("AG").ColumnWidth = 24.28515625
$ws.Columns("AH").ColumnWidth = 14.5703125
$ws.Columns("AI").ColumnWidth = 21.42578125
$ws.Columns("AJ").ColumnWidth = 16.85546875

# Update the view state to match the edited file.
$ws.Range("AE14").Select()
$excel.ActiveWindow.ScrollColumn = 25
